# Manutencao nos cenarios existentes
# Atualiza os valores de NOVO_EMAIL (I2) e NOVO_CPF (J2) para o cenario de
# cadastro de pessoa fisica.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# NOVO_EMAIL: texto simples, sem alteracao de estilo.
$ws.Range("I2").Value = "alexsantos_26032022110939@gmail.com"

# NOVO_CPF: precisa permanecer como texto (nao numero), mantendo o mesmo
# estilo da celula (sem formato de numero/texto aplicado). Usamos o prefixo
# de apostrofo para forcar o tipo texto e, em seguida, colamos apenas a
# formatacao original da celula (identica a de I2) para anular a mudanca de
# estilo causada pelo "quote prefix".
$ws.Range("J2").Value = "'85519994943"

$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
